$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "factura" rows to the new "servicio" wording
$ws.Range("D9").Value = "Registrar servicio"
$ws.Range("D10").Value = "Ver servicio"
$ws.Range("D11").Value = "Editar servicio"

# Add the new row 12 entry, matching formatting of the row above it
$ws.Range("D11:F11").Copy()
$null = $ws.Range("D12:F12").PasteSpecial(-4122)
$ws.Range("D12").Value = "Eliminar servicio"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""

$ws.Range("D12:F12").RowHeight = 15.75

# Update the selection shown in the sheet view
$null = $ws.Range("F18").Select()
